$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# "mais dados ciclo-2": E7 and E8 (cycle-2 rows) get their "done" hours
# entered as 0 (typed over the old shared formula that copied C7/C8), so
# G (the remaining balance, C-E) now actually computes instead of sitting
# at the stale 0 it had before.
$ws.Range("E7").Formula = "=0"
$ws.Range("E8").Value = 0

# Cursor/selection ends up on E9 in the saved file.
$ws.Range("E9").Select()
